$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 420, shifting existing rows 420:545 down to 421:546
$ws.Rows("420:420").Insert()

# Populate the newly inserted row 420 with its values
$ws.Range("A420").Value = 6
$ws.Range("B420").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C420").Value = "Metropolitana"
$ws.Range("D420").Value = 44841
$ws.Range("E420").Value = 13
$ws.Range("F420").Value = 100112039
$ws.Range("G420").Value = "Ciboulette"
$ws.Range("H420").Value = "Sin especificar"
$ws.Range("I420").Value = "Primera"
$ws.Range("J420").Value = 770
$ws.Range("K420").Value = 800
$ws.Range("L420").Value = 900
$ws.Range("M420").Value = 845
$ws.Range("N420").Value = "$/docena de atados"
$ws.Range("O420").Value = "Región Metropolitana"
$ws.Range("P420").Value = 282
$ws.Range("Q420").Value = 3
$ws.Range("R420").Value = "Hortaliza"
